$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Maintainer of the project" -> "Maintainer of the projects"
# ------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Maintainer of the project ", $true, $false, $false, $false, $false, $true, 1, $false, "Maintainer of the projects ", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Insert the new "northern-lights-forecast" hyperlink, followed by
#    " and ", right before the existing "ncdump-rich" hyperlink.
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("ncdump-rich", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.InsertBefore("NLF-PLACEHOLDER and ")

$r3 = $d.Content
$r3.Find.Execute("NLF-PLACEHOLDER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$h = $d.Hyperlinks.Add($r3, "https://github.com/engeir/northern-lights-forecast", "", "", "northern-lights-forecast", "")
$h.Range.Bold = 1
$h.Range.BoldBi = 1

# ------------------------------------------------------------------
# 3. "which is published on" -> "which are both published on"
#    (search starts inside the existing plain run, not right at the
#    hyperlink boundary, to avoid picking up Hyperlink character
#    formatting for the replacement text)
# ------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("which is published on", $true, $false, $false, $false, $false, $true, 1, $false, "which are both published on", 2) | Out-Null

# ------------------------------------------------------------------
# 4. ". This is a previewer for quickly showing formatted metadata in"
#    -> new, longer sentence describing both projects.
#    Again start the search just after the leading period so we never
#    anchor the match immediately after a hyperlink's closing tag.
# ------------------------------------------------------------------
$r5 = $d.Content
$r5.Find.Execute("This is a previewer for quickly showing formatted metadata in", $true, $false, $false, $false, $false, $true, 1, $false, "northern-lights-forecast is a program that listens to a website for updates on northern lights weather, and sends a message to a Telegram bot if conditions for seeing northern lights are good. ncdump-rich is a previewer for quickly showing formatted metadata in", 2) | Out-Null

Write-Output "done"
